# Auto-generated Excel COM-interop script
# Applies the scheduled market-price refresh described in the commit diff:
# updates currentAveragePrice/LevePrice/LeveProfit columns (H-N) for specific
# leve rows across all eight job sheets, matching the upstream OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 99999
$ws.Range("I12").Value = 99999
$ws.Range("K12").Value = 99999
$ws.Range("M12").Value = -99829
$ws.Range("H33").Value = 201.73334
$ws.Range("J33").Value = 500
$ws.Range("L33").Value = 500
$ws.Range("N33").Value = -958
$ws.Range("H40").Value = 5557482
$ws.Range("I40").Value = 2036.4615
$ws.Range("J40").Value = 15874737
$ws.Range("K40").Value = 2036.4615
$ws.Range("L40").Value = 15874737
$ws.Range("M40").Value = -1861.4615
$ws.Range("N40").Value = -15875087
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("H51").Value = 12189.875
$ws.Range("J51").Value = 12068.087
$ws.Range("L51").Value = 12068.087
$ws.Range("N51").Value = -13036.087
$ws.Range("H116").Value = 4539.6665
$ws.Range("J116").Value = 4377.7144
$ws.Range("L116").Value = 4377.7144
$ws.Range("N116").Value = -11261.7144
$ws.Range("H125").Value = 11366653
$ws.Range("I125").Value = 3118.7144
$ws.Range("J125").Value = 31252838
$ws.Range("K125").Value = 28068.4296
$ws.Range("L125").Value = 281275542
$ws.Range("M125").Value = -25608.4296
$ws.Range("N125").Value = -281280462
$ws.Range("H132").Value = 2272.9167
$ws.Range("I132").Value = 2294.8333
$ws.Range("K132").Value = 6884.499899999999
$ws.Range("M132").Value = -4354.499899999999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("L133").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2963.394
$ws.Range("I32").Value = 2912.4482
$ws.Range("J32").Value = 3332.75
$ws.Range("K32").Value = 2912.4482
$ws.Range("L32").Value = 3332.75
$ws.Range("M32").Value = -2625.4482
$ws.Range("N32").Value = -3906.75
$ws.Range("H39").Value = 21248.75
$ws.Range("I39").Value = 21248.75
$ws.Range("K39").Value = 21248.75
$ws.Range("M39").Value = -20728.75
$ws.Range("H63").Value = 2831.1667
$ws.Range("I63").Value = 2831.1667
$ws.Range("K63").Value = 2831.1667
$ws.Range("M63").Value = -2145.1667
$ws.Range("H66").Value = 2831.1667
$ws.Range("I66").Value = 2831.1667
$ws.Range("K66").Value = 14155.8335
$ws.Range("M66").Value = -10723.8335
$ws.Range("H102").Value = 1396.4286
$ws.Range("I102").Value = 1396.4286
$ws.Range("K102").Value = 1396.4286
$ws.Range("M102").Value = 225.5714
$ws.Range("H122").Value = 1573.1515
$ws.Range("I122").Value = 1276.08
$ws.Range("K122").Value = 3828.24
$ws.Range("M122").Value = -1378.24

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("H82").Value = 18000
$ws.Range("I82").Value = 18000
$ws.Range("K82").Value = 18000
$ws.Range("M82").Value = -17617
$ws.Range("H85").Value = 18000
$ws.Range("I85").Value = 18000
$ws.Range("K85").Value = 18000
$ws.Range("M85").Value = -16674
$ws.Range("H86").Value = 2701.5386
$ws.Range("I86").Value = 2641.111
$ws.Range("K86").Value = 2641.111
$ws.Range("M86").Value = -1518.111
$ws.Range("H89").Value = 2701.5386
$ws.Range("I89").Value = 2641.111
$ws.Range("K89").Value = 13205.555
$ws.Range("M89").Value = -7589.555
$ws.Range("H105").Value = 3105.4
$ws.Range("I105").Value = 3105.4
$ws.Range("K105").Value = 3105.4
$ws.Range("M105").Value = -1358.4
$ws.Range("M33").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1007.8
$ws.Range("I35").Value = 933.8
$ws.Range("K35").Value = 933.8
$ws.Range("M35").Value = -639.8
$ws.Range("H99").Value = 2039
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 2398.3333
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 2398.3333
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -5394.3333
$ws.Range("H126").Value = 2039
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2398.3333
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 7194.999899999999
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -12134.9999
$ws.Range("H132").Value = 100003570
$ws.Range("I132").Value = 100003570
$ws.Range("K132").Value = 300010710
$ws.Range("M132").Value = -300008180
$ws.Range("H133").Value = 93331.836
$ws.Range("J133").Value = 93331.836
$ws.Range("L133").Value = 93331.836
$ws.Range("N133").Value = -98391.836
$ws.Range("H134").Value = 19233868
$ws.Range("I134").Value = 22730298
$ws.Range("K134").Value = 68190894
$ws.Range("M134").Value = -68188359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1586.5714
$ws.Range("I34").Value = 550.5
$ws.Range("K34").Value = 1651.5
$ws.Range("M34").Value = -1567.5
$ws.Range("H55").Value = 948.5
$ws.Range("J55").Value = 1131.3334
$ws.Range("L55").Value = 3394.0002
$ws.Range("N55").Value = -3748.0002
$ws.Range("H56").Value = 16318.185
$ws.Range("I56").Value = 16318.185
$ws.Range("K56").Value = 16318.185
$ws.Range("M56").Value = -15788.185
$ws.Range("H86").Value = 707
$ws.Range("I86").Value = 519.75
$ws.Range("K86").Value = 1559.25
$ws.Range("M86").Value = -373.25
$ws.Range("H89").Value = 707
$ws.Range("I89").Value = 519.75
$ws.Range("K89").Value = 4677.75
$ws.Range("M89").Value = 1250.25
$ws.Range("H118").Value = 1558.4
$ws.Range("I118").Value = 1558.4
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 4675.200000000001
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = -3432.200000000001
$ws.Range("H121").Value = 103094.75
$ws.Range("I121").Value = 204257.8
$ws.Range("J121").Value = 30835.428
$ws.Range("K121").Value = 612773.3999999999
$ws.Range("L121").Value = 92506.284
$ws.Range("M121").Value = -611463.3999999999
$ws.Range("N121").Value = -95126.284
$ws.Range("H131").Value = 2048.2856
$ws.Range("J131").Value = 2819.8
$ws.Range("L131").Value = 8459.400000000001
$ws.Range("N131").Value = -18539.4
$ws.Range("H140").Value = 1795.8846
$ws.Range("I140").Value = 767.9545000000001
$ws.Range("K140").Value = 2303.8635
$ws.Range("M140").Value = 2876.1365
$ws.Range("M118").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("H132").Value = 11366750
$ws.Range("I132").Value = 11366750
$ws.Range("K132").Value = 34100250
$ws.Range("M132").Value = -34097720
$ws.Range("H140").Value = 89996.664
$ws.Range("J140").Value = 89996.664
$ws.Range("L140").Value = 89996.664
$ws.Range("N140").Value = -100356.664
$ws.Range("L38").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1724.8667
$ws.Range("I82").Value = 1765.3846
$ws.Range("K82").Value = 1765.3846
$ws.Range("M82").Value = -1404.3846
$ws.Range("H85").Value = 1724.8667
$ws.Range("I85").Value = 1765.3846
$ws.Range("K85").Value = 1765.3846
$ws.Range("M85").Value = -517.3846000000001
$ws.Range("H122").Value = 670791.1
$ws.Range("J122").Value = 3998.1667
$ws.Range("L122").Value = 11994.5001
$ws.Range("N122").Value = -16894.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 645.6923
$ws.Range("I113").Value = 480.57144
$ws.Range("K113").Value = 1441.71432
$ws.Range("M113").Value = 728.28568
$ws.Range("H126").Value = 1179.8928
$ws.Range("I126").Value = 1205.1923
$ws.Range("K126").Value = 3615.5769
$ws.Range("M126").Value = -1145.5769
$ws.Range("H132").Value = 14290728
$ws.Range("I132").Value = 20002174
$ws.Range("J132").Value = 12110.7
$ws.Range("K132").Value = 60006522
$ws.Range("L132").Value = 36332.10000000001
$ws.Range("M132").Value = -60003992
$ws.Range("N132").Value = -41392.10000000001

Write-Output "Applied scheduled Sheets price refresh."